$d = $word.ActiveDocument

# Step 1: remove the old _GoBack bookmark at its original location, and merge the
# two runs that it used to split into one (matches diff hunk near "curatescience.org").
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()
$oldText = "large scale replication curation projects like curatescience.org "
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $oldText, 2) | Out-Null

# Step 2: merge the 4 runs in the "See Addressing..." paragraph into a single run.
$seeText = "See Addressing the “Replication Crisis”: Using Original Studies to Design Replication Studies with Appropriate Statistical Power Samantha F. Anderson & Scott E. Maxwell"
$d.Content.Find.Execute($seeText, $true, $false, $false, $false, $false, $true, 1, $false, $seeText, 2) | Out-Null

# Step 3: insert the 5 new paragraphs (to-do notes) at the very start of the document.
$r = $d.Range(0, 0)
$r.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t>ALSO write another paper about the trend in effect size reporting – i.e., descriptive and over time.</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t xml:space='preserve'>Also write another paper about the trend in effect sizes reported in the literature over time. Are they in fact going </w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>down.</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r><w:t xml:space='preserve'> </w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t xml:space='preserve'>Maybe run them all through </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>statcheck</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t xml:space='preserve'> too? </w:t></w:r><w:r><w:t>Or at least the APA formatted ones.</w:t></w:r><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>")
